$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-07 Wednesday", "2026-01-08 Thursday"),
    @("479÷7=", "426÷9="),
    @("220÷8=", "881÷4="),
    @("589÷5=", "834÷8="),
    @("705÷3=", "636÷8="),
    @("458÷8=", "490÷9="),
    @("394÷6=", "619÷7="),
    @("725÷2=", "141÷6="),
    @("923÷9=", "963÷2="),
    @("339÷3=", "527÷5="),
    @("816÷7=", "676÷2="),
    @("309÷6=", "805÷4="),
    @("345÷9=", "900÷8="),
    @("676÷9=", "940÷5="),
    @("849÷5=", "328÷8="),
    @("958÷3=", "762÷2="),
    @("660÷3=", "442÷4="),
    @("823÷3=", "289÷7="),
    @("813÷5=", "709÷3="),
    @("825÷5=", "162÷4="),
    @("281÷2=", "492÷7="),
    @("387÷6=", "579÷5="),
    @("500÷2=", "985÷8="),
    @("771÷9=", "183÷4="),
    @("980÷5=", "707÷5="),
    @("183÷6=", "805÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
